# Insert a new row at position 20 (shifts existing rows 20-50 down to 21-51)
# and populate it with the new weekly record described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20, 1).Value = 3
$ws.Cells.Item(20, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 44571
$ws.Cells.Item(20, 5).Value = 5
$ws.Cells.Item(20, 6).Value = 100112022
$ws.Cells.Item(20, 7).Value = "Arveja Verde"
$ws.Cells.Item(20, 8).Value = "Perfection"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 73
$ws.Cells.Item(20, 11).Value = 15000
$ws.Cells.Item(20, 12).Value = 16000
$ws.Cells.Item(20, 13).Value = 15479
$ws.Cells.Item(20, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 16).Value = 619
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"
